$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.044705284417568
$ws.Cells.Item(2, 4).Value = 1.047349785034542
$ws.Cells.Item(2, 5).Value = 1.051737923372595
$ws.Cells.Item(2, 6).Value = 1.0606697192026
$ws.Cells.Item(2, 9).Value = 1.037779805521828
$ws.Cells.Item(2, 10).Value = 1.049768939149611
$ws.Cells.Item(2, 11).Value = 1.050112669623495
$ws.Cells.Item(2, 12).Value = 1.054488594488906
$ws.Cells.Item(2, 13).Value = 1.063395868606503
$ws.Cells.Item(2, 14).Value = 1.020367742962661

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.045952768263674
$ws.Cells.Item(3, 4).Value = 1.048512554600585
$ws.Cells.Item(3, 5).Value = 1.052849805471249
$ws.Cells.Item(3, 6).Value = 1.061873592691342
$ws.Cells.Item(3, 9).Value = 1.03800658802673
$ws.Cells.Item(3, 10).Value = 1.050662313420207
$ws.Cells.Item(3, 11).Value = 1.051086355557237
$ws.Cells.Item(3, 12).Value = 1.055412413816581
$ws.Cells.Item(3, 13).Value = 1.064413230265951
$ws.Cells.Item(3, 14).Value = 1.020672052285571

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.046759573573921
$ws.Cells.Item(4, 4).Value = 1.049264830246337
$ws.Cells.Item(4, 5).Value = 1.053569226102874
$ws.Cells.Item(4, 6).Value = 1.062652600757053
$ws.Cells.Item(4, 9).Value = 1.038151695725416
$ws.Cells.Item(4, 10).Value = 1.051239486644071
$ws.Cells.Item(4, 11).Value = 1.051715705986755
$ws.Cells.Item(4, 12).Value = 1.056009557178154
$ws.Cells.Item(4, 13).Value = 1.065070981112206
$ws.Cells.Item(4, 14).Value = 1.020868467069257

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.047098661206075
$ws.Cells.Item(5, 4).Value = 1.049581061338056
$ws.Cells.Item(5, 5).Value = 1.053871662095667
$ws.Cells.Item(5, 6).Value = 1.062980102542939
$ws.Cells.Item(5, 9).Value = 1.038212307745903
$ws.Cells.Item(5, 10).Value = 1.051481916329845
$ws.Cells.Item(5, 11).Value = 1.051980121371048
$ws.Cells.Item(5, 12).Value = 1.056260446915264
$ws.Cells.Item(5, 13).Value = 1.065347369334207
$ws.Cells.Item(5, 14).Value = 1.020950922033975

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.04715559011118
$ws.Cells.Item(6, 4).Value = 1.049634156478746
$ws.Cells.Item(6, 5).Value = 1.053922441977581
$ws.Cells.Item(6, 6).Value = 1.06303509198073
$ws.Cells.Item(6, 9).Value = 1.038222461837395
$ws.Cells.Item(6, 10).Value = 1.051522608837338
$ws.Cells.Item(6, 11).Value = 1.052024508331987
$ws.Cells.Item(6, 12).Value = 1.056302563689093
$ws.Cells.Item(6, 13).Value = 1.065393768538804
$ws.Cells.Item(6, 14).Value = 1.020964759687967

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.046764104843052
$ws.Cells.Item(7, 4).Value = 1.049269055840089
$ws.Cells.Item(7, 5).Value = 1.053573267298738
$ws.Cells.Item(7, 6).Value = 1.062656976823169
$ws.Cells.Item(7, 9).Value = 1.038152507163142
$ws.Cells.Item(7, 10).Value = 1.051242726840233
$ws.Cells.Item(7, 11).Value = 1.051719239758752
$ws.Cells.Item(7, 12).Value = 1.056012910164166
$ws.Cells.Item(7, 13).Value = 1.065074674735048
$ws.Cells.Item(7, 14).Value = 1.020869569298972

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.045126961771371
$ws.Cells.Item(8, 4).Value = 1.047742773089542
$ws.Cells.Item(8, 5).Value = 1.05211369860509
$ws.Cells.Item(8, 6).Value = 1.061076570455282
$ws.Cells.Item(8, 9).Value = 1.037856786631745
$ws.Cells.Item(8, 10).Value = 1.050071046320335
$ws.Cells.Item(8, 11).Value = 1.050441875024809
$ws.Cells.Item(8, 12).Value = 1.054800934731333
$ws.Cells.Item(8, 13).Value = 1.063739805552426
$ws.Cells.Item(8, 14).Value = 1.020470688176874

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.04223893053545
$ws.Cells.Item(9, 4).Value = 1.045052317762241
$ws.Cells.Item(9, 5).Value = 1.049541367973866
$ws.Cells.Item(9, 6).Value = 1.058291781356711
$ws.Cells.Item(9, 9).Value = 1.037323145653888
$ws.Cells.Item(9, 10).Value = 1.047999440004797
$ws.Cells.Item(9, 11).Value = 1.048185645756952
$ws.Cells.Item(9, 12).Value = 1.052660395735892
$ws.Cells.Item(9, 13).Value = 1.061383315046018
$ws.Cells.Item(9, 14).Value = 1.019764010202407

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.040311283217446
$ws.Cells.Item(10, 4).Value = 1.043257931933236
$ws.Cells.Item(10, 5).Value = 1.047826126637655
$ws.Cells.Item(10, 6).Value = 1.056435207164796
$ws.Cells.Item(10, 9).Value = 1.036958928076332
$ws.Cells.Item(10, 10).Value = 1.046613599201845
$ws.Cells.Item(10, 11).Value = 1.046677797632379
$ws.Cells.Item(10, 12).Value = 1.051229997239372
$ws.Cells.Item(10, 13).Value = 1.059809347558341
$ws.Cells.Item(10, 14).Value = 1.019290312385488

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.039476010916649
$ws.Cells.Item(11, 4).Value = 1.042480740089707
$ws.Cells.Item(11, 5).Value = 1.047083301994256
$ws.Cells.Item(11, 6).Value = 1.05563125377692
$ws.Cells.Item(11, 9).Value = 1.036799206373501
$ws.Cells.Item(11, 10).Value = 1.046012361958761
$ws.Cells.Item(11, 11).Value = 1.046023983569242
$ws.Cells.Item(11, 12).Value = 1.050609798410408
$ws.Cells.Item(11, 13).Value = 1.059127074780243
$ws.Cells.Item(11, 14).Value = 1.019084578304429

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.039165662027384
$ws.Cells.Item(12, 4).Value = 1.042192022461044
$ws.Cells.Item(12, 5).Value = 1.046807364736684
$ws.Cells.Item(12, 6).Value = 1.055332620762268
$ws.Cells.Item(12, 9).Value = 1.036739575630699
$ws.Cells.Item(12, 10).Value = 1.045788859595704
$ws.Cells.Item(12, 11).Value = 1.045780989763311
$ws.Cells.Item(12, 12).Value = 1.050379302860395
$ws.Cells.Item(12, 13).Value = 1.058873535939454
$ws.Cells.Item(12, 14).Value = 1.019008065833512

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.03923223711241
$ws.Cells.Item(13, 4).Value = 1.042253954946598
$ws.Cells.Item(13, 5).Value = 1.046866555112301
$ws.Cells.Item(13, 6).Value = 1.055396678997888
$ws.Cells.Item(13, 9).Value = 1.036752380347239
$ws.Cells.Item(13, 10).Value = 1.045836809610707
$ws.Cells.Item(13, 11).Value = 1.045833119044291
$ws.Cells.Item(13, 12).Value = 1.050428750674246
$ws.Cells.Item(13, 13).Value = 1.058927925997989
$ws.Cells.Item(13, 14).Value = 1.019024482266098

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.039450359267013
$ws.Cells.Item(14, 4).Value = 1.042456875306742
$ws.Cells.Item(14, 5).Value = 1.047060493335978
$ws.Cells.Item(14, 6).Value = 1.055606568873955
$ws.Cells.Item(14, 9).Value = 1.036794283461794
$ws.Cells.Item(14, 10).Value = 1.045993890787013
$ws.Cells.Item(14, 11).Value = 1.04600390045561
$ws.Cells.Item(14, 12).Value = 1.050590748159379
$ws.Cells.Item(14, 13).Value = 1.059106119484176
$ws.Cells.Item(14, 14).Value = 1.019078255673836

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.039584739296625
$ws.Cells.Item(15, 4).Value = 1.042581896667579
$ws.Cells.Item(15, 5).Value = 1.047179982483716
$ws.Cells.Item(15, 6).Value = 1.055735887703216
$ws.Cells.Item(15, 9).Value = 1.036820061188017
$ws.Cells.Item(15, 10).Value = 1.046090650293173
$ws.Cells.Item(15, 11).Value = 1.04610910615156
$ws.Cells.Item(15, 12).Value = 1.050690543392339
$ws.Cells.Item(15, 13).Value = 1.059215895424113
$ws.Cells.Item(15, 14).Value = 1.01911137481586

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.04036670482315
$ws.Cells.Item(16, 4).Value = 1.043309506957457
$ws.Cells.Item(16, 5).Value = 1.047875422893305
$ws.Cells.Item(16, 6).Value = 1.05648856177436
$ws.Cells.Item(16, 9).Value = 1.036969485819644
$ws.Cells.Item(16, 10).Value = 1.046653476751887
$ws.Cells.Item(16, 11).Value = 1.04672116985758
$ws.Cells.Item(16, 12).Value = 1.051271140213319
$ws.Cells.Item(16, 13).Value = 1.059854612110078
$ws.Cells.Item(16, 14).Value = 1.019303953173119

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.040857051392501
$ws.Cells.Item(17, 4).Value = 1.043765859586135
$ws.Cells.Item(17, 5).Value = 1.048311622398175
$ws.Cells.Item(17, 6).Value = 1.056960680835818
$ws.Cells.Item(17, 9).Value = 1.03706267659704
$ws.Cells.Item(17, 10).Value = 1.047006211240243
$ws.Cells.Item(17, 11).Value = 1.047104857273438
$ws.Cells.Item(17, 12).Value = 1.051635110575273
$ws.Cells.Item(17, 13).Value = 1.06025506415956
$ws.Cells.Item(17, 14).Value = 1.019424586159856

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.041143005705391
$ws.Cells.Item(18, 4).Value = 1.044032022156945
$ws.Cells.Item(18, 5).Value = 1.048566039344319
$ws.Cells.Item(18, 6).Value = 1.057236055755311
$ws.Cells.Item(18, 9).Value = 1.037116839004586
$ws.Cells.Item(18, 10).Value = 1.047211843603037
$ws.Cells.Item(18, 11).Value = 1.047328568388133
$ws.Cells.Item(18, 12).Value = 1.051847328678949
$ws.Cells.Item(18, 13).Value = 1.060488570336425
$ws.Cells.Item(18, 14).Value = 1.019494889576399

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.041240499240998
$ws.Cells.Item(19, 4).Value = 1.044122773358059
$ws.Cells.Item(19, 5).Value = 1.048652787168896
$ws.Cells.Item(19, 6).Value = 1.05732995085129
$ws.Cells.Item(19, 9).Value = 1.037135274063142
$ws.Cells.Item(19, 10).Value = 1.047281940037527
$ws.Cells.Item(19, 11).Value = 1.047404833357215
$ws.Cells.Item(19, 12).Value = 1.051919676113926
$ws.Cells.Item(19, 13).Value = 1.060568177971757
$ws.Cells.Item(19, 14).Value = 1.019518851104256

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.040804447719221
$ws.Cells.Item(20, 4).Value = 1.04371689937462
$ws.Cells.Item(20, 5).Value = 1.048264823456809
$ws.Cells.Item(20, 6).Value = 1.05691002737435
$ws.Cells.Item(20, 9).Value = 1.037052698199526
$ws.Cells.Item(20, 10).Value = 1.046968377736234
$ws.Cells.Item(20, 11).Value = 1.047063700286838
$ws.Cells.Item(20, 12).Value = 1.051596068253101
$ws.Cells.Item(20, 13).Value = 1.060212106775816
$ws.Cells.Item(20, 14).Value = 1.019411649568018

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.03938613028232
$ws.Cells.Item(21, 4).Value = 1.042397121274552
$ws.Cells.Item(21, 5).Value = 1.047003383906626
$ws.Cells.Item(21, 6).Value = 1.055544761801406
$ws.Cells.Item(21, 9).Value = 1.036781952407402
$ws.Cells.Item(21, 10).Value = 1.045947639172132
$ws.Cells.Item(21, 11).Value = 1.045953613416802
$ws.Cells.Item(21, 12).Value = 1.050543047433444
$ws.Cells.Item(21, 13).Value = 1.059053649061501
$ws.Cells.Item(21, 14).Value = 1.01906242333583

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.03849384645965
$ws.Cells.Item(22, 4).Value = 1.041567126703526
$ws.Cells.Item(22, 5).Value = 1.046210155224069
$ws.Cells.Item(22, 6).Value = 1.054686312362295
$ws.Cells.Item(22, 9).Value = 1.036609970497776
$ws.Cells.Item(22, 10).Value = 1.045304840602523
$ws.Cells.Item(22, 11).Value = 1.045254857414827
$ws.Cells.Item(22, 12).Value = 1.049880240997905
$ws.Cells.Item(22, 13).Value = 1.058324629678385
$ws.Cells.Item(22, 14).Value = 1.01884230870292

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.038966914428819
$ws.Cells.Item(23, 4).Value = 1.042007142034176
$ws.Cells.Item(23, 5).Value = 1.046630672006587
$ws.Cells.Item(23, 6).Value = 1.055141398407775
$ws.Cells.Item(23, 9).Value = 1.03670130773674
$ws.Cells.Item(23, 10).Value = 1.045645697597802
$ws.Cells.Item(23, 11).Value = 1.04562535774283
$ws.Cells.Item(23, 12).Value = 1.050231677094473
$ws.Cells.Item(23, 13).Value = 1.058711159030719
$ws.Cells.Item(23, 14).Value = 1.018959047218971

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.040828217229534
$ws.Cells.Item(24, 4).Value = 1.043739022449163
$ws.Cells.Item(24, 5).Value = 1.048285969917527
$ws.Cells.Item(24, 6).Value = 1.056932915506067
$ws.Cells.Item(24, 9).Value = 1.037057207607908
$ws.Cells.Item(24, 10).Value = 1.046985473413303
$ws.Cells.Item(24, 11).Value = 1.047082297627143
$ws.Cells.Item(24, 12).Value = 1.051613710044062
$ws.Cells.Item(24, 13).Value = 1.060231517587017
$ws.Cells.Item(24, 14).Value = 1.019417495241701

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.042985948423521
$ws.Cells.Item(25, 4).Value = 1.045747988856171
$ws.Cells.Item(25, 5).Value = 1.050206432472267
$ws.Cells.Item(25, 6).Value = 1.059011716635626
$ws.Cells.Item(25, 9).Value = 1.037462593509058
$ws.Cells.Item(25, 10).Value = 1.048535833750497
$ws.Cells.Item(25, 11).Value = 1.048769579163023
$ws.Cells.Item(25, 12).Value = 1.053214364238099
$ws.Cells.Item(25, 13).Value = 1.061993041425344
$ws.Cells.Item(25, 14).Value = 1.019947156127129

